$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Correct mis-shifted input values in columns Q/T (rows 53-61) ---
$ws.Range("T53").Value = 0.19020000000000001
$ws.Range("T54").Value = 28.88
$ws.Range("T55").Value = 383.1
$ws.Range("T56").Value = 928.9
$ws.Range("T57").Value = 1552
$ws.Range("T58").Value = 3454

$ws.Range("Q59").Value = 12830
$ws.Range("T59").Value = 5162

$ws.Range("Q60").Value = 12860
$ws.Range("T60").Value = 6515

$ws.Range("Q61").Value = 12560
$ws.Range("T61").Value = 7461

# --- Switch column AH (rows 40-48 and 53-61) from scientific to fixed 2-decimal format ---
$ws.Range("AH40:AH48").NumberFormat = "0.00"
$ws.Range("AH53:AH61").NumberFormat = "0.00"

# --- Update the active selection on the sheet ---
$ws.Range("I43").Select()
